# g13.3a: round "Valor" column to 2 decimals, swap Nordeste/Brasil row order,
# add thin border + top-vertical alignment to header row, reset page margins
# to Excel defaults.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round the "Valor" (D) column values down to 2 decimal places ---
$ws.Range("D2").Value = 97.82
$ws.Range("D3").Value = 97.72
$ws.Range("D4").Value = 97.24
$ws.Range("D5").Value = 96.53
$ws.Range("D6").Value = 95.98
$ws.Range("D7").Value = 95.89
$ws.Range("D8").Value = 91.61

# --- Swap the order of the "Nordeste" / "Brasil" summary rows (9 & 10) ---
$ws.Range("A9").Value = "Brasil"
$ws.Range("D9").Value = 93.64

$ws.Range("A10").Value = "Nordeste"
$ws.Range("D10").Value = 91.33

# --- Header row formatting: thin box border + vertical-top alignment ---
$header = $ws.Range("A1:E1")
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# --- Reset page margins to Excel's standard defaults (values are in points) ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
